# Weekly update: insert a new price-report row for "Acelga" (chard) at
# Macroferia Regional de Talca, pushing the existing rows 118:169 down to
# 119:170 (dimension grows from A1:R169 to A1:R170).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row above the current row 118; Excel shifts rows
# 118:169 down to 119:170 and copies row 118's formatting into the new row.
$ws.Rows.Item(118).Insert()

# Fill the new row 118 with the latest week's record.
$ws.Cells.Item(118, 1).Value  = 5
$ws.Cells.Item(118, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(118, 3).Value  = "Maule"
$ws.Cells.Item(118, 4).Value  = 44460
$ws.Cells.Item(118, 5).Value  = 7
$ws.Cells.Item(118, 6).Value  = 100112009
$ws.Cells.Item(118, 7).Value  = "Acelga"
$ws.Cells.Item(118, 8).Value  = "Sin especificar"
$ws.Cells.Item(118, 9).Value  = "Primera"
$ws.Cells.Item(118, 10).Value = 500
$ws.Cells.Item(118, 11).Value = 2300
$ws.Cells.Item(118, 12).Value = 2300
$ws.Cells.Item(118, 13).Value = 2300
$ws.Cells.Item(118, 14).Value = "`$/docena de atados (4 kilos)"
$ws.Cells.Item(118, 15).Value = "Región del Maule"
$ws.Cells.Item(118, 16).Value = 575
$ws.Cells.Item(118, 17).Value = 4
$ws.Cells.Item(118, 18).Value = "Hortaliza"

# Keep the date cell formatted like the rest of column D.
$ws.Cells.Item(118, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
